$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column N (Grand Total shifts from N to O),
# which also shifts the header cell style/formatting along with it.
$ws.Columns("N").Insert()

# Fill in the header for the newly inserted column with the "PPN" label.
$ws.Range("N7").Value = "PPN"

# Match the author's final selection after adding the column.
$ws.Range("O9").Select() | Out-Null
